$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Shift the last two columns (NOME DO MENTOR, EMAIL DO MENTOR) one column to the
# right (E->F, F->G) to make room for the new "Nivel" column at D, working from
# the rightmost column first so we don't clobber data before it is copied.
for ($r = 1; $r -le 8; $r++) {
    $ws.Cells.Item($r, 7).Value2 = $ws.Cells.Item($r, 6).Value2   # G <- F (EMAIL DO MENTOR)
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 5).Value2   # F <- E (NOME DO MENTOR)
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 4).Value2   # E <- D (CIDADE)
}

# New column D: header + value for every data row.
$ws.Range("D1").Value2 = "Nível"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "Superior"
}

# E5 (Pindamonhangaba, moved from D5) gets underlined.
$ws.Range("E5").Font.Underline = 2

# Update the active selection to match the edited cell.
$ws.Range("E5").Select() | Out-Null
